$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2021 column (R) to the table, matching the existing
# formatting used for the preceding year column (Q).

# Header cell R4: year 2021, using the same style as Q4.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("R4").Value = 2021

# Data cell R5: value 3.6, using the same style as Q5.
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("R5").Value = 3.6

# Update the selected cell shown in the saved view.
$ws.Range("O9").Select() | Out-Null
